$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0127ea0209efe6467dc7a75a3f6a35c8a08cf1a/e2e/7c960b42-9426-4e3f-b83f-7187c679ef85.md"
$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0127ea0209efe6467dc7a75a3f6a35c8a08cf1a/e2e/7c960b42-9426-4e3f-b83f-7187c679ef85.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40a59b00c273a441ff5f1c6ce55708081de2c325/e2e/7c960b42-9426-4e3f-b83f-7187c679ef85.md."
$displayMd = "7c960b42-9426-4e3f-b83f-7187c679ef85.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(16).ColumnWidth = 39.17

$wsZh.Range("J8").Value = "7c960b42-9426-4e3f-b83f-7187c679ef85.d08d55193343ab78c42b9962e000676759a9c02f.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-16 04:41:43"
$wsZh.Range("P8").Value = $errorMsg

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $handbackUrl, "", "", $displayMd)
$wsZh.Range("I8").Font.Name = "Calibri"
$wsZh.Range("I8").Font.Underline = $true
$wsZh.Range("I8").Font.Color = 15570276

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(16).ColumnWidth = 39.17

$wsDe.Range("J8").Value = "7c960b42-9426-4e3f-b83f-7187c679ef85.d08d55193343ab78c42b9962e000676759a9c02f.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-16 04:41:50"
$wsDe.Range("P8").Value = $errorMsg

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $handbackUrl, "", "", $displayMd)
$wsDe.Range("I8").Font.Name = "Calibri"
$wsDe.Range("I8").Font.Underline = $true
$wsDe.Range("I8").Font.Color = 15570276
